# Charts with Sqr. Mileage
# Apply the data edits: fix "North West" -> "NW" label, fill in the two
# missing region totals (J6/J8), drop the "Multiple"/"META-ANALYSIS" rows
# (I12:J12, I13:J13), and add the new "Regions2"/"Sq. mi2 (mi^2)" columns
# (K:L) mirroring the Regions/Sq. mi (mi^2) data in I:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "North West" region label (I5) ---
$ws.Range("I5").Value = "NW"

# --- Fill in the two region totals that were "N/A" ---
$ws.Range("J6").Value = 855767
$ws.Range("J8").Value = 895300

# --- Remove the "Multiple" / "META-ANALYSIS" helper rows ---
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()

# --- Add the new Regions2 / Sq. mi2 (mi^2) columns (K:L) ---
$ws.Range("K1").Value = "Regions2"
$ws.Range("L1").Value = "Sq. mi2 (mi^2)"

$ws.Range("K2").Value = "East-North Central"
$ws.Range("L2").Value = 299170

$ws.Range("K3").Value = "Mid-Atlantic"
$ws.Range("L3").Value = 191308

$ws.Range("K4").Value = "Mountain"
$ws.Range("L4").Value = 855767

$ws.Range("K5").Value = "New England"
$ws.Range("L5").Value = 71992

$ws.Range("K6").Value = "NW"
$ws.Range("L6").Value = 5469

$ws.Range("K7").Value = "Pacific"
$ws.Range("L7").Value = 895300

$ws.Range("K8").Value = "South Atlantic"
$ws.Range("L8").Value = 292589

$ws.Range("K9").Value = "West North Central"
$ws.Range("L9").Value = 507900

$ws.Range("K10").Value = "West South Central"
$ws.Range("L10").Value = 444100

# --- Match the saved selection (top-left viewport cell isn't exposed via
#     this interop surface, but the active-cell selection is) ---
$ws.Range("L1").Select()
